$wb = $excel.ActiveWorkbook

# Sheet ALC, row 2 (Leve Item ID 5489)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 493.44446
$ws.Range("I2").Value = 323.66666
$ws.Range("K2").Value = 323.66666
$ws.Range("M2").Value = -210.66666

# Sheet ALC, row 17 (Leve Item ID 38956)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4345.1724
$ws.Range("J17").Value = 4345.1724
$ws.Range("L17").Value = 13035.5172
$ws.Range("N17").Value = -13371.5172

# Sheet ALC, row 33 (Leve Item ID 5512)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 285.77777
$ws.Range("I33").Value = 340
$ws.Range("K33").Value = 340
$ws.Range("M33").Value = -111

# Sheet ALC, row 58 (Leve Item ID 4606)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 112785
$ws.Range("J58").Value = 202980
$ws.Range("L58").Value = 608940
$ws.Range("N58").Value = -609240

# Sheet ALC, row 64 (Leve Item ID 5506)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3753.8462
$ws.Range("I64").Value = 3300
$ws.Range("J64").Value = 4037.5
$ws.Range("K64").Value = 3300
$ws.Range("L64").Value = 4037.5
$ws.Range("M64").Value = -3052
$ws.Range("N64").Value = -4533.5

# Sheet ALC, row 67 (Leve Item ID 5506)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3753.8462
$ws.Range("I67").Value = 3300
$ws.Range("J67").Value = 4037.5
$ws.Range("K67").Value = 3300
$ws.Range("L67").Value = 4037.5
$ws.Range("M67").Value = -2442
$ws.Range("N67").Value = -5753.5

# Sheet ALC, row 129 (Leve Item ID 36115)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 644.6667
$ws.Range("I129").Value = 412.2857
$ws.Range("J129").Value = 1458
$ws.Range("K129").Value = 1236.8571
$ws.Range("L129").Value = 4374
$ws.Range("M129").Value = 3763.1429
$ws.Range("N129").Value = -14374

# Sheet ALC, row 138 (Leve Item ID 44169)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1073573.8
$ws.Range("J138").Value = 1361547.4
$ws.Range("L138").Value = 4084642.2
$ws.Range("N138").Value = -4094922.2

# Sheet ARM, row 32 (Leve Item ID 44147)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11591.378
$ws.Range("I32").Value = 8670.825000000001
$ws.Range("J32").Value = 28318.182
$ws.Range("K32").Value = 8670.825000000001
$ws.Range("L32").Value = 28318.182
$ws.Range("M32").Value = -8383.825000000001
$ws.Range("N32").Value = -28892.182

# Sheet BSM, row 94 (Leve Item ID 19939)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1531.591
$ws.Range("I94").Value = 1315.8334
$ws.Range("J94").Value = 2502.5
$ws.Range("K94").Value = 1315.8334
$ws.Range("L94").Value = 2502.5
$ws.Range("M94").Value = -864.8334
$ws.Range("N94").Value = -3404.5

# Sheet BSM, row 102 (Leve Item ID 19565)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H102").Value = 17028
$ws.Range("I102").Value = 4750
$ws.Range("J102").Value = 29306
$ws.Range("K102").Value = 4750
$ws.Range("L102").Value = 29306
$ws.Range("M102").Value = -1505
$ws.Range("N102").Value = -35796

# Sheet BSM, row 134 (Leve Item ID 43998)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 19450.648
$ws.Range("I134").Value = 1857.6909
$ws.Range("J134").Value = 503257
$ws.Range("K134").Value = 5573.072700000001
$ws.Range("L134").Value = 1509771
$ws.Range("M134").Value = -3038.072700000001
$ws.Range("N134").Value = -1514841

# Sheet BSM, row 140 (Leve Item ID 42471)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 56390.562
$ws.Range("J140").Value = 56390.562
$ws.Range("L140").Value = 56390.562
$ws.Range("N140").Value = -66750.56200000001

# Sheet CRP, row 16 (Leve Item ID 27691)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1450.4615
$ws.Range("I16").Value = 1428
$ws.Range("J16").Value = 1501
$ws.Range("K16").Value = 1428
$ws.Range("L16").Value = 1501
$ws.Range("M16").Value = -1141
$ws.Range("N16").Value = -2075

# Sheet CRP, row 31 (Leve Item ID 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 14610.6
$ws.Range("I31").Value = 3525
$ws.Range("J31").Value = 22001
$ws.Range("K31").Value = 3525
$ws.Range("L31").Value = 22001
$ws.Range("M31").Value = -3230
$ws.Range("N31").Value = -22591

# Sheet CRP, row 34 (Leve Item ID 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 14610.6
$ws.Range("I34").Value = 3525
$ws.Range("J34").Value = 22001
$ws.Range("K34").Value = 3525
$ws.Range("L34").Value = 22001
$ws.Range("M34").Value = -3323
$ws.Range("N34").Value = -22405

# Sheet CRP, row 113 (Leve Item ID 27691)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1450.4615
$ws.Range("I113").Value = 1428
$ws.Range("J113").Value = 1501
$ws.Range("K113").Value = 1428
$ws.Range("L113").Value = 1501
$ws.Range("M113").Value = 742
$ws.Range("N113").Value = -5841

# Sheet CUL, row 105 (Leve Item ID 19814)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H105").Value = 5309.3335
$ws.Range("J105").Value = 5309.3335
$ws.Range("L105").Value = 15928.0005
$ws.Range("N105").Value = -21170.0005

# Sheet GSM, row 21 (Leve Item ID 4430)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 1501.75
$ws.Range("I21").Value = 1000
$ws.Range("J21").Value = 3007
$ws.Range("K21").Value = 1000
$ws.Range("L21").Value = 3007
$ws.Range("M21").Value = -827
$ws.Range("N21").Value = -3353

# Sheet GSM, row 29 (Leve Item ID 4209)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 4429.5713
$ws.Range("I29").Value = 2201.4
$ws.Range("K29").Value = 2201.4
$ws.Range("M29").Value = -1911.4

# Sheet GSM, row 30 (Leve Item ID 4430)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H30").Value = 1501.75
$ws.Range("I30").Value = 1000
$ws.Range("J30").Value = 3007
$ws.Range("K30").Value = 1000
$ws.Range("L30").Value = 3007
$ws.Range("M30").Value = -895
$ws.Range("N30").Value = -3217

# Sheet GSM, row 113 (Leve Item ID 27710)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2824.9524
$ws.Range("I113").Value = 2800.6875
$ws.Range("J113").Value = 2902.6
$ws.Range("K113").Value = 2800.6875
$ws.Range("L113").Value = 2902.6
$ws.Range("M113").Value = -630.6875
$ws.Range("N113").Value = -7242.6

# Sheet LTW, row 16 (Leve Item ID 5289)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 704.06665
$ws.Range("I16").Value = 704.06665
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 704.06665
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -534.06665

# Sheet LTW, row 21 (Leve Item ID 2672)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H21").Value = 11000
$ws.Range("J21").Value = 11000
$ws.Range("L21").Value = 11000
$ws.Range("N21").Value = -11348

# Sheet LTW, row 23 (Leve Item ID 4097)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 403196
$ws.Range("I23").Value = 403196
$ws.Range("K23").Value = 403196
$ws.Range("M23").Value = -402966

# Sheet LTW, row 31 (Leve Item ID 3043)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H31").Value = 517.8570999999999
$ws.Range("I31").Value = 517.8570999999999
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 517.8570999999999
$ws.Range("L31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -269.8570999999999

# Sheet LTW, row 40 (Leve Item ID 36248)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4226.8
$ws.Range("I40").Value = 3912.16
$ws.Range("K40").Value = 3912.16
$ws.Range("M40").Value = -3776.16

# Sheet LTW, row 46 (Leve Item ID 5282)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1084
$ws.Range("J46").Value = 1034.2858
$ws.Range("L46").Value = 1034.2858
$ws.Range("N46").Value = -1410.2858

# Sheet LTW, row 61 (Leve Item ID 27740)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1569110.1
$ws.Range("I61").Value = 1569110.1
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1569110.1
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -1568908.1

# Sheet LTW, row 93 (Leve Item ID 19993)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 9500
$ws.Range("I93").Value = 10000
$ws.Range("K93").Value = 10000
$ws.Range("M93").Value = -8752

# Sheet LTW, row 113 (Leve Item ID 27740)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1569110.1
$ws.Range("I113").Value = 1569110.1
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1569110.1
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -1566940.1

# Sheet LTW, row 122 (Leve Item ID 36247)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5656.95
$ws.Range("I122").Value = 4561.3105
$ws.Range("J122").Value = 8545.454
$ws.Range("K122").Value = 13683.9315
$ws.Range("L122").Value = 25636.362
$ws.Range("M122").Value = -11233.9315
$ws.Range("N122").Value = -30536.362

# Sheet WVR, row 23 (Leve Item ID 3325)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 5920.1665
$ws.Range("I23").Value = 2170
$ws.Range("J23").Value = 9670.333000000001
$ws.Range("K23").Value = 2170
$ws.Range("L23").Value = 9670.333000000001
$ws.Range("M23").Value = -1941
$ws.Range("N23").Value = -10128.333

# Sheet WVR, row 132 (Leve Item ID 44029)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2919.353
$ws.Range("I132").Value = 1941.3636
$ws.Range("J132").Value = 4712.3335
$ws.Range("K132").Value = 5824.0908
$ws.Range("L132").Value = 14137.0005
$ws.Range("M132").Value = -3294.0908
$ws.Range("N132").Value = -19197.0005
